$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.082.88"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "3.141.35"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D5").Value = "'589.27"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "'137.80"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.137.65"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'5.24"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "'34.16"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "3.656.94"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "3.138.67"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "63.059.75"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "'6.66"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").Value = "'472.25"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'84.81"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").Value = "'13.04"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'2.71"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'7.05"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'7.96"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'26.81"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'0.106"
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("D37").Value = "'52.27"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "0.0₃0690"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'420.21"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "  -6.67%  "
$ws.Range("D42").Value = "'8.19"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "2.912.50"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D45").Value = "'0.262"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "'25.43"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("E50").Value = "  -6.95%  "
$ws.Range("D51").Value = "'120.29"
$ws.Range("E51").Value = "  -0.90%  "
